$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 363.4643
$ws.Range("J96").Value = 498.8889
$ws.Range("L96").Value = 1496.6667
$ws.Range("N96").Value = -4242.6667

$ws.Range("H103").Value = 8696233
$ws.Range("I103").Value = 458
$ws.Range("J103").Value = 15385290
$ws.Range("K103").Value = 1374
$ws.Range("L103").Value = 46155870
$ws.Range("M103").Value = -788
$ws.Range("N103").Value = -46157042

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5239.12
$ws.Range("I32").Value = 4184.6597
$ws.Range("J32").Value = 39333.332
$ws.Range("K32").Value = 4184.6597
$ws.Range("L32").Value = 39333.332
$ws.Range("M32").Value = -3897.6597
$ws.Range("N32").Value = -39907.332

$ws.Range("H61").Value = 1068.0625
$ws.Range("I61").Value = 1021.871
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 1021.871
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -809.871
$ws.Range("N61").Value = -2924

$ws.Range("H74").Value = 2188.625
$ws.Range("I74").Value = 834.6667
$ws.Range("J74").Value = 9500
$ws.Range("K74").Value = 834.6667
$ws.Range("L74").Value = 9500
$ws.Range("M74").Value = 39.33330000000001
$ws.Range("N74").Value = -11248

$ws.Range("H77").Value = 2188.625
$ws.Range("I77").Value = 834.6667
$ws.Range("J77").Value = 9500
$ws.Range("K77").Value = 4173.3335
$ws.Range("L77").Value = 47500
$ws.Range("M77").Value = 194.6665000000003
$ws.Range("N77").Value = -56236

$ws.Range("H102").Value = 50001590
$ws.Range("I102").Value = 1538.7858
$ws.Range("J102").Value = 166668380
$ws.Range("K102").Value = 1538.7858
$ws.Range("L102").Value = 166668380
$ws.Range("M102").Value = 83.21419999999989
$ws.Range("N102").Value = -166671624

$ws.Range("H122").Value = 14216.4
$ws.Range("I122").Value = 14216.4
$ws.Range("K122").Value = 42649.2
$ws.Range("M122").Value = -40199.2

$ws.Range("H132").Value = 28574184
$ws.Range("I132").Value = 33334736
$ws.Range("J132").Value = 10868.4
$ws.Range("K132").Value = 100004208
$ws.Range("L132").Value = 32605.2
$ws.Range("M132").Value = -100001678
$ws.Range("N132").Value = -37665.2

$ws.Range("H136").Value = 1068.0625
$ws.Range("I136").Value = 1021.871
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 3065.613
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -515.6129999999998
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 44058.832
$ws.Range("J74").Value = 44058.832
$ws.Range("L74").Value = 44058.832
$ws.Range("N74").Value = -45930.832

$ws.Range("H77").Value = 44058.832
$ws.Range("J77").Value = 44058.832
$ws.Range("L77").Value = 132176.496
$ws.Range("N77").Value = -141536.496

$ws.Range("H94").Value = 553.4737
$ws.Range("I94").Value = 502.66666
$ws.Range("J94").Value = 678.1818
$ws.Range("K94").Value = 502.66666
$ws.Range("L94").Value = 678.1818
$ws.Range("M94").Value = -51.66665999999998
$ws.Range("N94").Value = -1580.1818

$ws.Range("H103").Value = 42399.8
$ws.Range("J103").Value = 42399.8
$ws.Range("L103").Value = 42399.8
$ws.Range("N103").Value = -44743.8

$ws.Range("H128").Value = 1400
$ws.Range("I128").Value = 1400
$ws.Range("K128").Value = 4200
$ws.Range("M128").Value = -1710

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2313.8462
$ws.Range("I122").Value = 1466.6666
$ws.Range("J122").Value = 3040
$ws.Range("K122").Value = 4399.9998
$ws.Range("L122").Value = 9120
$ws.Range("M122").Value = -1949.9998
$ws.Range("N122").Value = -14020

$ws.Range("H134").Value = 6217.6875
$ws.Range("I134").Value = 6777.357
$ws.Range("J134").Value = 2300
$ws.Range("K134").Value = 20332.071
$ws.Range("L134").Value = 6900
$ws.Range("M134").Value = -17797.071
$ws.Range("N134").Value = -11970

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1474.326
$ws.Range("I122").Value = 1358.3948
$ws.Range("J122").Value = 2025
$ws.Range("K122").Value = 4075.1844
$ws.Range("L122").Value = 6075
$ws.Range("M122").Value = -1625.1844
$ws.Range("N122").Value = -10975

$ws.Range("H132").Value = 2652.258
$ws.Range("I132").Value = 2462.35
$ws.Range("J132").Value = 2997.5454
$ws.Range("K132").Value = 7387.049999999999
$ws.Range("L132").Value = 8992.636200000001
$ws.Range("M132").Value = -4857.049999999999
$ws.Range("N132").Value = -14052.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7695418.5
$ws.Range("I7").Value = 11113059
$ws.Range("K7").Value = 11113059
$ws.Range("M7").Value = -11112947

$ws.Range("H100").Value = 2250
$ws.Range("I100").Value = 1800
$ws.Range("J100").Value = 3600
$ws.Range("K100").Value = 1800
$ws.Range("L100").Value = 3600
$ws.Range("M100").Value = -1259
$ws.Range("N100").Value = -4682

$ws.Range("H122").Value = 3726.56
$ws.Range("I122").Value = 3264.5715
$ws.Range("J122").Value = 4314.5454
$ws.Range("K122").Value = 9793.7145
$ws.Range("L122").Value = 12943.6362
$ws.Range("M122").Value = -7343.7145
$ws.Range("N122").Value = -17843.6362

$ws.Range("H126").Value = 7695418.5
$ws.Range("I126").Value = 11113059
$ws.Range("K126").Value = 33339177
$ws.Range("M126").Value = -33336707

$ws.Range("H132").Value = 7476.385
$ws.Range("I132").Value = 8188.4443
$ws.Range("J132").Value = 5874.25
$ws.Range("K132").Value = 24565.3329
$ws.Range("L132").Value = 17622.75
$ws.Range("M132").Value = -22035.3329
$ws.Range("N132").Value = -22682.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 30000
$ws.Range("J86").Value = 30000
$ws.Range("L86").Value = 30000
$ws.Range("N86").Value = -32246

$ws.Range("H89").Value = 30000
$ws.Range("J89").Value = 30000
$ws.Range("L89").Value = 150000
$ws.Range("N89").Value = -161232

$ws.Range("H96").Value = 1031.8695
$ws.Range("I96").Value = 993.93335
$ws.Range("J96").Value = 1103
$ws.Range("K96").Value = 993.93335
$ws.Range("L96").Value = 1103
$ws.Range("M96").Value = 379.06665
$ws.Range("N96").Value = -3849

$ws.Range("H122").Value = 9641.4
$ws.Range("I122").Value = 12658.647
$ws.Range("J122").Value = 3229.75
$ws.Range("K122").Value = 37975.94100000001
$ws.Range("L122").Value = 9689.25
$ws.Range("M122").Value = -35525.94100000001
$ws.Range("N122").Value = -14589.25
